# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the data block (rows 43-44),
# pushing the existing rows 43-54 down to rows 45-56.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 43, shifting rows 43:54 down to 45:56
$ws.Rows("43:44").Insert()

# New row 43
$ws.Cells.Item(43,1).Value = 10
$ws.Cells.Item(43,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(43,3).Value = "La Araucanía"
$ws.Cells.Item(43,4).Value = 44508
$ws.Cells.Item(43,5).Value = 9
$ws.Cells.Item(43,6).Value = 100112022
$ws.Cells.Item(43,7).Value = "Arveja Verde"
$ws.Cells.Item(43,8).Value = "Sin especificar"
$ws.Cells.Item(43,9).Value = "Primera"
$ws.Cells.Item(43,10).Value = 30
$ws.Cells.Item(43,11).Value = 14000
$ws.Cells.Item(43,12).Value = 14000
$ws.Cells.Item(43,13).Value = 14000
$ws.Cells.Item(43,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(43,15).Value = "Región Metropolitana"
$ws.Cells.Item(43,16).Value = 560
$ws.Cells.Item(43,17).Value = 25
$ws.Cells.Item(43,18).Value = "Hortaliza"

# New row 44
$ws.Cells.Item(44,1).Value = 10
$ws.Cells.Item(44,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(44,3).Value = "La Araucanía"
$ws.Cells.Item(44,4).Value = 44508
$ws.Cells.Item(44,5).Value = 9
$ws.Cells.Item(44,6).Value = 100112022
$ws.Cells.Item(44,7).Value = "Arveja Verde"
$ws.Cells.Item(44,8).Value = "Sin especificar"
$ws.Cells.Item(44,9).Value = "Primera"
$ws.Cells.Item(44,10).Value = 70
$ws.Cells.Item(44,11).Value = 16000
$ws.Cells.Item(44,12).Value = 17000
$ws.Cells.Item(44,13).Value = 16571
$ws.Cells.Item(44,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(44,15).Value = "Región del Maule"
$ws.Cells.Item(44,16).Value = 663
$ws.Cells.Item(44,17).Value = 25
$ws.Cells.Item(44,18).Value = "Hortaliza"

Write-Host "Rows inserted and populated."
